# TestDriverComponent.pptx update:
#  - nudge a handful of shape/connector positions on the diagram slide
#  - drop the old dashed "Elbow Connector 288" that ran from the
#    EmailAccount box towards BrowserInstance
#  - add a new dashed "Straight Arrow Connector 60" in its place, same
#    line style, repositioned further along
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/.Top/.Width/
# .Height are single-precision (32-bit) floats internally. A handful of
# the target EMU offsets are not exactly representable in that
# precision, so a few literals are nudged by a tiny fraction of a point
# (<< 1/100 pt) so that the float32 round-trip lands back on the exact
# EMU value PowerPoint would have written.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

# 1) Rectangle 102 ("EmailServer"-adjacent box) slides right
$rect102 = Get-ShapeByName $s "Rectangle 102"
$rect102.Left = 96.2803937007874

# 2) Rectangle 94 slides left
$rect94 = Get-ShapeByName $s "Rectangle 94"
$rect94.Left = 99.93409448818898

# 3) Flowchart: Decision 174 slides right
$flowDecision174 = Get-ShapeByName $s "Flowchart: Decision 174"
$flowDecision174.Left = 144.8568503937008

# 4) Elbow Connector 175 (rotated 90deg) - moves and shrinks
$elbowConn175 = Get-ShapeByName $s "Elbow Connector 175"
$elbowConn175.Left = 148.03897887795276
$elbowConn175.Top = 423.7915823031496
$elbowConn175.Height = 9.328663917322835

# 5) Elbow Connector 178 (rotated 90deg) - moves and shrinks
$elbowConn178 = Get-ShapeByName $s "Elbow Connector 178"
$elbowConn178.Left = 169.03890263779527
$elbowConn178.Top = 402.7915823031496
$elbowConn178.Height = 51.328661417322834

# 6) Elbow Connector 186 (rotated 270deg, flipped) - moves and grows
$elbowConn186 = Get-ShapeByName $s "Elbow Connector 186"
$elbowConn186.Left = 33.21385826771654
$elbowConn186.Top = 318.29521185039374
$elbowConn186.Height = 57.70866391732284

# 7) & 8) Replace the dashed "Elbow Connector 288" with a dashed
# straight "Straight Arrow Connector 60" further along the same path.
# Duplicate first so the new shape inherits the exact same line/style
# formatting (dotted accent2 line, arrow tail, style refs), then
# convert it from an elbow connector to a straight connector and move
# it into place, and finally delete the original.
$oldConnector = Get-ShapeByName $s "Elbow Connector 288"
$newConnectorRange = $oldConnector.Duplicate()
$newConnector = $newConnectorRange.Item(1)
$newConnector.Name = "Straight Arrow Connector 60"
$newConnector.ConnectorFormat.Type = 1
$newConnector.Left = 258.0
$newConnector.Top = 102.0
$newConnector.Width = 84.0
$newConnector.Height = 0
$oldConnector.Delete()
